$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 17544078
$ws.Range("I6").Value = 33333434
$ws.Range("K6").Value = 100000302
$ws.Range("M6").Value = -100000190

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 10666.667
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 15000
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = -1831
$ws.Range("N13").Value = -15338

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 44532.332
$ws.Range("I62").Value = 29444
$ws.Range("K62").Value = 29444
$ws.Range("M62").Value = -28820

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 44532.332
$ws.Range("I65").Value = 29444
$ws.Range("K65").Value = 147220
$ws.Range("M65").Value = -144100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7145726
$ws.Range("I70").Value = 18185182
$ws.Range("J70").Value = 2548.9412
$ws.Range("K70").Value = 54555546
$ws.Range("L70").Value = 7646.823600000001
$ws.Range("M70").Value = -54555276
$ws.Range("N70").Value = -8186.823600000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 7145726
$ws.Range("I73").Value = 18185182
$ws.Range("J73").Value = 2548.9412
$ws.Range("K73").Value = 54555546
$ws.Range("L73").Value = 7646.823600000001
$ws.Range("M73").Value = -54554610
$ws.Range("N73").Value = -9518.8236

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 7174.25
$ws.Range("I106").Value = 8161.625
$ws.Range("J106").Value = 6186.875
$ws.Range("K106").Value = 8161.625
$ws.Range("L106").Value = 6186.875
$ws.Range("M106").Value = -7530.625
$ws.Range("N106").Value = -7448.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2018442.9
$ws.Range("J116").Value = 3996.6667
$ws.Range("L116").Value = 3996.6667
$ws.Range("N116").Value = -10880.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3661147.8
$ws.Range("I137").Value = 543394.75
$ws.Range("J137").Value = 6952109.5
$ws.Range("K137").Value = 1630184.25
$ws.Range("L137").Value = 20856328.5
$ws.Range("M137").Value = -1627634.25
$ws.Range("N137").Value = -20861428.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 9447.518
$ws.Range("I138").Value = 8416.583000000001
$ws.Range("K138").Value = 25249.749
$ws.Range("M138").Value = -20109.749

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7311.4375
$ws.Range("I141").Value = 5221.6665
$ws.Range("J141").Value = 9998.286
$ws.Range("K141").Value = 15664.9995
$ws.Range("L141").Value = 29994.858
$ws.Range("M141").Value = -10484.9995
$ws.Range("N141").Value = -40354.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2087.762
$ws.Range("I32").Value = 2010.6
$ws.Range("K32").Value = 2010.6
$ws.Range("M32").Value = -1723.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12331.143
$ws.Range("I61").Value = 18805.285
$ws.Range("K61").Value = 18805.285
$ws.Range("M61").Value = -18593.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 55717284
$ws.Range("J74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("N74").Value = -5248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 55717284
$ws.Range("J77").Value = 3500
$ws.Range("L77").Value = 17500
$ws.Range("N77").Value = -26236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2332.7778
$ws.Range("I110").Value = 1222
$ws.Range("J110").Value = 3443.5557
$ws.Range("K110").Value = 1222
$ws.Range("L110").Value = 3443.5557
$ws.Range("M110").Value = 823
$ws.Range("N110").Value = -7533.5557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4073.4443
$ws.Range("I132").Value = 3410.2778
$ws.Range("J132").Value = 5399.778
$ws.Range("K132").Value = 10230.8334
$ws.Range("L132").Value = 16199.334
$ws.Range("M132").Value = -7700.8334
$ws.Range("N132").Value = -21259.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 12331.143
$ws.Range("I136").Value = 18805.285
$ws.Range("K136").Value = 56415.855
$ws.Range("M136").Value = -53865.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 12165.923
$ws.Range("I99").Value = 13805.552
$ws.Range("K99").Value = 13805.552
$ws.Range("M99").Value = -12307.552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9638.833000000001
$ws.Range("I105").Value = 12136.546
$ws.Range("J105").Value = 5713.857
$ws.Range("K105").Value = 12136.546
$ws.Range("L105").Value = 5713.857
$ws.Range("M105").Value = -10389.546
$ws.Range("N105").Value = -9207.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2922.1853
$ws.Range("I107").Value = 2991.2173
$ws.Range("J107").Value = 2525.25
$ws.Range("K107").Value = 2991.2173
$ws.Range("L107").Value = 2525.25
$ws.Range("M107").Value = -1071.2173
$ws.Range("N107").Value = -6365.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 5416.6665
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5416.6665
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5416.6665
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -5764.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2493.0815
$ws.Range("I31").Value = 1781.4783
$ws.Range("J31").Value = 3122.577
$ws.Range("K31").Value = 1781.4783
$ws.Range("L31").Value = 3122.577
$ws.Range("M31").Value = -1486.4783
$ws.Range("N31").Value = -3712.577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2493.0815
$ws.Range("I34").Value = 1781.4783
$ws.Range("J34").Value = 3122.577
$ws.Range("K34").Value = 1781.4783
$ws.Range("L34").Value = 3122.577
$ws.Range("M34").Value = -1579.4783
$ws.Range("N34").Value = -3526.577

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 250938.75
$ws.Range("J5").Value = 556944
$ws.Range("L5").Value = 1670832
$ws.Range("N5").Value = -1671056

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 125000150
$ws.Range("I7").Value = 174.2
$ws.Range("J7").Value = 333333440
$ws.Range("K7").Value = 522.5999999999999
$ws.Range("L7").Value = 1000000320
$ws.Range("M7").Value = -410.5999999999999
$ws.Range("N7").Value = -1000000544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 16672530
$ws.Range("J68").Value = 23817272
$ws.Range("L68").Value = 71451816
$ws.Range("N68").Value = -71453438

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 16672530
$ws.Range("J71").Value = 23817272
$ws.Range("L71").Value = 214355448
$ws.Range("N71").Value = -214363560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2525.6316
$ws.Range("J86").Value = 2691.3845
$ws.Range("L86").Value = 8074.1535
$ws.Range("N86").Value = -10446.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 2525.6316
$ws.Range("J89").Value = 2691.3845
$ws.Range("L89").Value = 24222.4605
$ws.Range("N89").Value = -36078.4605

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 992.7105
$ws.Range("J113").Value = 1220.36
$ws.Range("L113").Value = 3661.08
$ws.Range("N113").Value = -8001.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3604.077
$ws.Range("I134").Value = 2036.1428
$ws.Range("K134").Value = 6108.428400000001
$ws.Range("M134").Value = -1038.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 250938.75
$ws.Range("J135").Value = 556944
$ws.Range("L135").Value = 5012496
$ws.Range("N135").Value = -5017566

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 24565624
$ws.Range("I11").Value = 30837500
$ws.Range("J11").Value = 5750000
$ws.Range("K11").Value = 30837500
$ws.Range("L11").Value = 5750000
$ws.Range("M11").Value = -30837361
$ws.Range("N11").Value = -5750278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10675.676
$ws.Range("I12").Value = 10000
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -9860
$ws.Range("N12").Value = -15280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 185277.5
$ws.Range("I20").Value = 628451.25
$ws.Range("J20").Value = 8008
$ws.Range("K20").Value = 628451.25
$ws.Range("L20").Value = 8008
$ws.Range("M20").Value = -628206.25
$ws.Range("N20").Value = -8498

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5267.074
$ws.Range("I132").Value = 4168.68
$ws.Range("K132").Value = 12506.04
$ws.Range("M132").Value = -9976.040000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 19504
$ws.Range("I3").Value = 19504
$ws.Range("K3").Value = 19504
$ws.Range("M3").Value = -19392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 19504
$ws.Range("I15").Value = 19504
$ws.Range("K15").Value = 19504
$ws.Range("M15").Value = -19334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 93884.664
$ws.Range("I40").Value = 123329.664
$ws.Range("K40").Value = 123329.664
$ws.Range("M40").Value = -123193.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 571.32355
$ws.Range("I55").Value = 565.5769
$ws.Range("J55").Value = 590
$ws.Range("K55").Value = 565.5769
$ws.Range("L55").Value = 590
$ws.Range("M55").Value = -392.5769
$ws.Range("N55").Value = -936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2371.5
$ws.Range("I93").Value = 2964.182
$ws.Range("J93").Value = 1647.1111
$ws.Range("K93").Value = 2964.182
$ws.Range("L93").Value = 1647.1111
$ws.Range("M93").Value = -1716.182
$ws.Range("N93").Value = -4143.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4031.8965
$ws.Range("I100").Value = 4071.2964
$ws.Range("K100").Value = 4071.2964
$ws.Range("M100").Value = -3530.2964

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20421.188
$ws.Range("I132").Value = 23228.25
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 69684.75
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -67154.75
$ws.Range("N132").Value = -41060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 1000
$ws.Range("I7").Value = 1000
$ws.Range("K7").Value = 1000
$ws.Range("M7").Value = -887

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 3253333.2
$ws.Range("I9").Value = 2505000
$ws.Range("K9").Value = 2505000
$ws.Range("M9").Value = -2504860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2520.1875
$ws.Range("I107").Value = 2529.889
$ws.Range("K107").Value = 7589.667
$ws.Range("M107").Value = -5669.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3039.6667
$ws.Range("I113").Value = 1015.2105
$ws.Range("J113").Value = 7847.75
$ws.Range("K113").Value = 3045.6315
$ws.Range("L113").Value = 23543.25
$ws.Range("M113").Value = -875.6315
$ws.Range("N113").Value = -27883.25
